# Apply updated crypto price/volume data and fix row order for rows 48-49
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Formula = "59.138.66"
$ws.Range("E2").Formula = "  +1.45%  "
$ws.Range("D3").Formula = "2.633.61"
$ws.Range("E3").Formula = "  +3.74%  "
$ws.Range("D4").Formula = "1.00"
$ws.Range("E4").Formula = "  +0.03%  "
$ws.Range("D5").Formula = "522.56"
$ws.Range("E5").Formula = "  +3.46%  "
$ws.Range("D6").Formula = "146.40"
$ws.Range("E6").Formula = "  +2.27%  "
$ws.Range("D7").Formula = "0.995"
$ws.Range("E7").Formula = "  -0.39%  "
$ws.Range("E8").Formula = "  +1.14%  "
$ws.Range("D9").Formula = "2.660.29"
$ws.Range("E9").Formula = "  +4.53%  "
$ws.Range("D10").Formula = "6.32"
$ws.Range("E10").Formula = "  +3.76%  "
$ws.Range("E11").Formula = "  +4.07%  "
$ws.Range("E12").Formula = "  +2.48%  "
$ws.Range("E13").Formula = "  -1.28%  "
$ws.Range("D14").Formula = "3.097.84"
$ws.Range("E14").Formula = "  +3.67%  "
$ws.Range("D15").Formula = "59.121.72"
$ws.Range("E15").Formula = "  +1.42%  "
$ws.Range("D16").Formula = "21.05"
$ws.Range("E16").Formula = "  +2.04%  "
$ws.Range("E17").Formula = "  +2.58%  "
$ws.Range("D18").Formula = "2.656.73"
$ws.Range("E18").Formula = "  +4.22%  "
$ws.Range("D19").Formula = "349.64"
$ws.Range("E19").Formula = "  +2.36%  "
$ws.Range("E20").Formula = "  +0.62%  "
$ws.Range("D21").Formula = "10.34"
$ws.Range("E21").Formula = "  +3.14%  "
$ws.Range("E22").Formula = "  +4.50%  "
$ws.Range("D23").Formula = "0.998"
$ws.Range("E23").Formula = "  +0.03%  "
$ws.Range("D24").Formula = "61.72"
$ws.Range("E24").Formula = "  +2.07%  "
$ws.Range("D25").Formula = "0.419"
$ws.Range("E25").Formula = "  +2.73%  "
$ws.Range("E26").Formula = "  +3.82%  "
$ws.Range("D27").Formula = "0.996"
$ws.Range("E27").Formula = "  -0.04%  "
$ws.Range("D28").Formula = "0.0₃0810"
$ws.Range("E28").Formula = "  +4.16%  "
$ws.Range("D29").Formula = "7.13"
$ws.Range("E29").Formula = "  +3.21%  "
$ws.Range("E30").Formula = "  -0.28%  "
$ws.Range("E31").Formula = "  +8.36%  "
$ws.Range("E32").Formula = "  +4.06%  "
$ws.Range("E33").Formula = "  +2.99%  "
$ws.Range("D34").Formula = "150.08"
$ws.Range("E34").Formula = "  +0.25%  "
$ws.Range("D35").Formula = "0.989"
$ws.Range("E35").Formula = "  +10.27%  "
$ws.Range("E36").Formula = "  +3.83%  "
$ws.Range("E37").Formula = "  +3.20%  "
$ws.Range("D38").Formula = "36.78"
$ws.Range("E38").Formula = "  +2.36%  "
$ws.Range("D39").Formula = "0.850"
$ws.Range("E39").Formula = "  +4.56%  "
$ws.Range("E40").Formula = "  +5.50%  "
$ws.Range("E41").Formula = "  +2.96%  "
$ws.Range("D42").Formula = "279.27"
$ws.Range("E42").Formula = "  -0.58%  "
$ws.Range("D43").Formula = "0.0987"
$ws.Range("E43").Formula = "  -0.14%  "
$ws.Range("D44").Formula = "0.994"
$ws.Range("E44").Formula = "  -0.40%  "
$ws.Range("D45").Formula = "0.609"
$ws.Range("E45").Formula = "  +1.77%  "
$ws.Range("D46").Formula = "19.66"
$ws.Range("E46").Formula = "  +5.57%  "
$ws.Range("D47").Formula = "0.0526"
$ws.Range("E47").Formula = "  -0.85%  "
$ws.Range("B48").Formula = "RenderToken"
$ws.Range("C48").Formula = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Formula = "4.77"
$ws.Range("E48").Formula = "  +5.97%  "
$ws.Range("B49").Formula = "VeChain"
$ws.Range("C49").Formula = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Formula = "0.0231"
$ws.Range("E49").Formula = "  +2.64%  "
$ws.Range("E50").Formula = "  +0.13%  "
$ws.Range("D51").Formula = "1.993.26"
$ws.Range("E51").Formula = "  +4.90%  "
